$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure new cells are treated as text (preserve leading zeros / exact string values)
$newRange = $ws.Range("A13:H20")
$newRange.NumberFormat = "@"

$ws.Range("A13").Value = "address"
$ws.Range("B13").Value = "HCM city"
$ws.Range("C13").Value = " "
$ws.Range("D13").Value = "first"
$ws.Range("E13").Value = "last"
$ws.Range("F13").Value = "039494944"
$ws.Range("G13").Value = "2"
$ws.Range("H13").Value = "{{credit_debit_cvv}}"
$ws.Range("A14").Value = "address"
$ws.Range("B14").Value = "H C M City "
$ws.Range("C14").Value = " "
$ws.Range("D14").Value = "firsst"
$ws.Range("E14").Value = "last"
$ws.Range("F14").Value = " "
$ws.Range("G14").Value = "222"
$ws.Range("H14").Value = " "
$ws.Range("A15").Value = "ddaddress"
$ws.Range("B15").Value = "dd city"
$ws.Range("C15").Value = " "
$ws.Range("D15").Value = "ddirst"
$ws.Range("E15").Value = "{{address}}"
$ws.Range("F15").Value = " "
$ws.Range("G15").Value = " "
$ws.Range("H15").Value = "{{address}}"
$ws.Range("A16").Value = "15 address"
$ws.Range("B16").Value = "15 HCM city"
$ws.Range("C16").Value = " "
$ws.Range("D16").Value = "15 silicon first"
$ws.Range("E16").Value = "15 silicon last"
$ws.Range("F16").Value = "{{ip_address}}"
$ws.Range("G16").Value = "1515"
$ws.Range("H16").Value = "151515"
$ws.Range("A17").Value = "15 only address"
$ws.Range("B17").Value = "15 only city"
$ws.Range("C17").Value = " "
$ws.Range("D17").Value = "15 only first"
$ws.Range("E17").Value = "15 only last"
$ws.Range("F17").Value = " "
$ws.Range("G17").Value = "{{address}}"
$ws.Range("H17").Value = "{{address}}"
$ws.Range("A18").Value = "f1 adrress"
$ws.Range("B18").Value = "01 city"
$ws.Range("C18").Value = " "
$ws.Range("D18").Value = "f1 first"
$ws.Range("E18").Value = "f1 last"
$ws.Range("F18").Value = "{{phone}}"
$ws.Range("G18").Value = "{{date_time}}"
$ws.Range("H18").Value = "{{date_time}}"
$ws.Range("A19").Value = "f03 address"
$ws.Range("B19").Value = "03 city"
$ws.Range("C19").Value = " "
$ws.Range("D19").Value = "f03 first"
$ws.Range("E19").Value = "f03 last"
$ws.Range("F19").Value = " "
$ws.Range("G19").Value = "0303"
$ws.Range("H19").Value = "030303"
$ws.Range("A20").Value = "address {{address}}"
$ws.Range("B20").Value = "f4 address"
$ws.Range("C20").Value = " "
$ws.Range("D20").Value = "first f4"
$ws.Range("E20").Value = "last f4"
$ws.Range("F20").Value = " "
$ws.Range("G20").Value = "{{pin}}"
$ws.Range("H20").Value = "{{in_aadhaar}}"

# Reset style to Normal so the appended rows do not carry the temporary text format
$newRange.Style = "Normal"
